$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price update reshuffles the D,M,N,O,P,R,S values across rows 2-5
# (row r's new values come from old row perm[r]).
# Capture the "before" values first, then write the "after" values.

$rows = 2,3,4,5

$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2   # D: Fecha
        M = $ws.Cells.Item($r, 13).Value2  # M: Volumen
        N = $ws.Cells.Item($r, 14).Value2  # N: Precio minimo
        O = $ws.Cells.Item($r, 15).Value2  # O: Precio maximo
        P = $ws.Cells.Item($r, 16).Value2  # P: Precio promedio ponderado
        R = $ws.Cells.Item($r, 18).Value2  # R: Origen
        S = $ws.Cells.Item($r, 19).Value2  # S: Precio $/Kg
    }
}

# permutation: new row r gets the old values of row perm[r]
$perm = @{ 2 = 3; 3 = 5; 4 = 2; 5 = 4 }

foreach ($r in $rows) {
    $src = $before[$perm[$r]]
    $ws.Cells.Item($r, 4).Value2 = $src.D
    $ws.Cells.Item($r, 13).Value2 = $src.M
    $ws.Cells.Item($r, 14).Value2 = $src.N
    $ws.Cells.Item($r, 15).Value2 = $src.O
    $ws.Cells.Item($r, 16).Value2 = $src.P
    $ws.Cells.Item($r, 18).Value2 = $src.R
    $ws.Cells.Item($r, 19).Value2 = $src.S
}
